$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing it to remain a Text cell and
# leaving the cell's style/number-format untouched (matches source data
# where numeric-looking price strings like "134.52" are stored as text).
function Set-TextValue($range, $value) {
    $savedStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $savedStyle
}

Set-TextValue $ws.Range("D2") '57.368.55'
$ws.Range("E2").Value = '  +1.68%  '
Set-TextValue $ws.Range("D3") '2.326.56'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue $ws.Range("D5") '541.92'
$ws.Range("E5").Value = '  +5.84%  '
Set-TextValue $ws.Range("D6") '134.52'
$ws.Range("E6").Value = '  +1.80%  '
$ws.Range("E7").Value = '  -0.62%  '
Set-TextValue $ws.Range("D8") '0.536'
$ws.Range("E8").Value = '  +0.71%  '
Set-TextValue $ws.Range("D9") '2.354.96'
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("E11").Value = '  +0.95%  '
Set-TextValue $ws.Range("D12") '5.38'
$ws.Range("E12").Value = '  +1.78%  '
Set-TextValue $ws.Range("D13") '0.355'
$ws.Range("E13").Value = '  +4.72%  '
Set-TextValue $ws.Range("D14") '2.751.30'
$ws.Range("E14").Value = '  +0.20%  '
Set-TextValue $ws.Range("D15") '23.52'
$ws.Range("E15").Value = '  -0.04%  '
Set-TextValue $ws.Range("D16") '57.435.48'
$ws.Range("E16").Value = '  +1.82%  '
Set-TextValue $ws.Range("D17") '0.0000134'
$ws.Range("E17").Value = '  +0.82%  '
Set-TextValue $ws.Range("D18") '2.342.25'
$ws.Range("E18").Value = '  +0.40%  '
Set-TextValue $ws.Range("D19") '335.44'
$ws.Range("E19").Value = '  +3.54%  '
Set-TextValue $ws.Range("D20") '10.55'
$ws.Range("E20").Value = '  +0.98%  '
Set-TextValue $ws.Range("D21") '4.21'
$ws.Range("E21").Value = '  +1.68%  '
Set-TextValue $ws.Range("D22") '6.78'
$ws.Range("E22").Value = '  +1.80%  '
Set-TextValue $ws.Range("D23") '0.997'
$ws.Range("E23").Value = '  -0.01%  '
Set-TextValue $ws.Range("D24") '61.73'
$ws.Range("E24").Value = '  -0.08%  '
Set-TextValue $ws.Range("D25") '0.170'
$ws.Range("E25").Value = '  +4.27%  '
$ws.Range("B26").Value = 'InternetComputer(DFINITY)'
$ws.Range("C26").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range("D26") '8.44'
$ws.Range("E26").Value = '  -2.58%  '
$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue $ws.Range("D27") '0.994'
$ws.Range("E27").Value = '  -0.41%  '
$ws.Range("E28").Value = '  +9.12%  '
$ws.Range("E29").Value = '  +4.80%  '
Set-TextValue $ws.Range("D30") '170.57'
$ws.Range("E30").Value = '  +1.84%  '
Set-TextValue $ws.Range("D31") '0.0₃0736'
$ws.Range("E31").Value = '  +2.50%  '
$ws.Range("E32").Value = '  +1.05%  '
Set-TextValue $ws.Range("D33") '18.57'
$ws.Range("E33").Value = '  +1.57%  '
Set-TextValue $ws.Range("D34") '1.02'
$ws.Range("E34").Value = '  +14.86%  '
Set-TextValue $ws.Range("D36") '0.990'
$ws.Range("E36").Value = '  -0.80%  '
Set-TextValue $ws.Range("D37") '1.26'
$ws.Range("E37").Value = '  -0.22%  '
Set-TextValue $ws.Range("D38") '4.12'
$ws.Range("E38").Value = '  +4.73%  '
$ws.Range("E39").Value = '  +3.97%  '
Set-TextValue $ws.Range("D40") '39.36'
Set-TextValue $ws.Range("D41") '148.77'
$ws.Range("E41").Value = '  -1.46%  '
$ws.Range("E42").Value = '  +0.87%  '
$ws.Range("E43").Value = '  +1.24%  '
Set-TextValue $ws.Range("D44") '281.63'
$ws.Range("E44").Value = '  +1.10%  '
Set-TextValue $ws.Range("D45") '19.33'
$ws.Range("E45").Value = '  +7.36%  '
Set-TextValue $ws.Range("D46") '0.0931'
$ws.Range("E46").Value = '  +0.66%  '
Set-TextValue $ws.Range("D47") '0.0506'
$ws.Range("E47").Value = '  +2.08%  '
$ws.Range("E48").Value = '  +0.32%  '
$ws.Range("E49").Value = '  +2.22%  '
$ws.Range("E50").Value = '  +1.21%  '
Set-TextValue $ws.Range("D51") '17.50'
$ws.Range("E51").Value = '  +2.06%  '
